$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect column D (price strings like "1.002") from COM auto-numeric coercion
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "22.448.21"
$ws.Range("E2").Value = "  +0.28%  "

$ws.Range("D3").Value = "1.573.42"
$ws.Range("E3").Value = "  +0.11%  "

$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("E5").Value = "  +0.00%  "

$ws.Range("D6").Value = "291.83"
$ws.Range("E6").Value = "  +0.33%  "

$ws.Range("D7").Value = "0.3731"
$ws.Range("E7").Value = "  -0.70%  "

$ws.Range("D8").Value = "49.97"
$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").Value = "0.3400"
$ws.Range("E9").Value = "  -0.54%  "

$ws.Range("B10").Value = "Polygon"
$ws.Range("C10").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D10").Value = "1.145"
$ws.Range("E10").Value = "  -0.44%  "

$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").Value = "0.07566"
$ws.Range("E11").Value = "  -0.98%  "

$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  -0.05%  "

$ws.Range("D13").Value = "21.32"
$ws.Range("E13").Value = "  +0.69%  "

$ws.Range("D14").Value = "6.021"
$ws.Range("E14").Value = "  -0.12%  "

$ws.Range("D15").Value = "6.969"
$ws.Range("E15").Value = "  +0.46%  "

$ws.Range("D16").Value = "1.572.53"
$ws.Range("E16").Value = "  -0.02%  "

$ws.Range("D17").Value = "0.00001121"
$ws.Range("E17").Value = "  -0.74%  "

$ws.Range("D18").Value = "90.96"
$ws.Range("E18").Value = "  +1.16%  "

$ws.Range("D19").Value = "0.06748"
$ws.Range("E19").Value = "  +0.23%  "

$ws.Range("E20").Value = "  +0.09%  "

$ws.Range("D21").Value = "6.306"
$ws.Range("E21").Value = "  +1.79%  "

$ws.Range("D22").Value = "16.33"
$ws.Range("E22").Value = "  -2.64%  "

$ws.Range("D23").Value = "12.18"
$ws.Range("E23").Value = "  +1.62%  "

$ws.Range("D24").Value = "22.464.77"
$ws.Range("E24").Value = "  +0.32%  "

$ws.Range("D25").Value = "2.340"
$ws.Range("E25").Value = "  -2.41%  "

$ws.Range("D26").Value = "2.693"
$ws.Range("E26").Value = "  +0.56%  "

$ws.Range("D27").Value = "20.13"
$ws.Range("E27").Value = "  -0.24%  "

$ws.Range("D28").Value = "148.60"
$ws.Range("E28").Value = "  +0.99%  "

$ws.Range("D29").Value = "5.008"
$ws.Range("E29").Value = "  -0.42%  "

$ws.Range("D30").Value = "125.57"
$ws.Range("E30").Value = "  -0.39%  "

$ws.Range("D31").Value = "1.748.01"
$ws.Range("E31").Value = "  +0.01%  "

$ws.Range("E32").Value = "  +7.53%  "

$ws.Range("D33").Value = "6.187"
$ws.Range("E33").Value = "  +0.80%  "

$ws.Range("D34").Value = "1.984"
$ws.Range("E34").Value = "  -1.14%  "

$ws.Range("D35").Value = "9.819"
$ws.Range("E35").Value = "  -0.53%  "

$ws.Range("D36").Value = "0.08380"
$ws.Range("E36").Value = "  -1.40%  "

$ws.Range("D37").Value = "1.374"
$ws.Range("E37").Value = "  +0.76%  "

$ws.Range("D38").Value = "0.02495"
$ws.Range("E38").Value = "  -1.75%  "

$ws.Range("D39").Value = "0.2299"
$ws.Range("E39").Value = "  -0.76%  "

$ws.Range("D40").Value = "0.06519"
$ws.Range("E40").Value = "  -0.66%  "

$ws.Range("D41").Value = "5.471"
$ws.Range("E41").Value = "  +1.23%  "

$ws.Range("D42").Value = "11.28"
$ws.Range("E42").Value = "  -1.28%  "

$ws.Range("E43").Value = "  -2.55%  "

$ws.Range("D44").Value = "1.002"
$ws.Range("E44").Value = "  +0.01%  "

$ws.Range("D45").Value = "14.05"
$ws.Range("E45").Value = "  -0.08%  "

$ws.Range("D46").Value = "3.810"
$ws.Range("E46").Value = "  +0.85%  "

$ws.Range("D47").Value = "0.5808"
$ws.Range("E47").Value = "  -2.85%  "

$ws.Range("D48").Value = "129.95"
$ws.Range("E48").Value = "  +3.61%  "

$ws.Range("E49").Value = "  -0.42%  "

$ws.Range("D50").Value = "1.226"
$ws.Range("E50").Value = "  -5.00%  "

$ws.Range("D51").Value = "0.07325"
$ws.Range("E51").Value = "  -0.01%  "

# Restore default style on column D so unaffected/affected cells keep original formatting
$ws.Range("D2:D51").Style = "Normal"
